$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# deletes
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("C4").ClearContents()

# updates
$ws.Range("C6").Value = 0.348613976222456
$ws.Range("C8").Value = -0.1384957661262676
$ws.Range("C9").Value = 1.891565607550105
$ws.Range("E9").Value = 1.609625625600009
$ws.Range("C10").Value = 1.566479473280191
$ws.Range("E10").Value = 2.49756057493542
$ws.Range("C11").Value = 1.113165545862116
$ws.Range("E11").Value = 1.609625625599986
$ws.Range("C12").Value = 0.7307568962937161
$ws.Range("E13").Value = 3.238605209600021
$ws.Range("C14").Value = 0.8188188121642126
$ws.Range("C15").Value = 1.384186838979828
$ws.Range("E15").Value = 2.777885851461526
$ws.Range("C17").Value = 2.349355943833076
$ws.Range("C18").Value = 1.9846842782967
$ws.Range("E18").Value = 2.047428048848809
$ws.Range("C19").Value = 1.78642563555842
$ws.Range("E20").Value = 1.552965246735782
$ws.Range("E21").Value = 1.216098605743343
$ws.Range("E22").Value = 0.232608152956959
$ws.Range("C23").Value = 0.5221702820068952
$ws.Range("C24").Value = 1.282262557986469
$ws.Range("E24").Value = 1.784618024189033
$ws.Range("C25").Value = 1.238324979098082
$ws.Range("E25").Value = 1.281608622679209
$ws.Range("C26").Value = 1.064321453542272
$ws.Range("E27").Value = 0.8660061896410554
$ws.Range("C28").Value = 2.247109253368307
$ws.Range("E29").Value = 0.5495555957892195
$ws.Range("C30").Value = 1.361817904277718
$ws.Range("C31").Value = 0.6311979695890368
$ws.Range("E31").Value = 0.07482640125564544
$ws.Range("C32").Value = -4.247034401476779
$ws.Range("E32").Value = -12.19860234240002
$ws.Range("C33").Value = -9.171727975571519
$ws.Range("C34").Value = -4.352425014431327
$ws.Range("E34").Value = 31.54369540926345
$ws.Range("C35").Value = -7.006249401853603
$ws.Range("E35").Value = -12.66856409363488
$ws.Range("C37").Value = -2.664090177971856
$ws.Range("E37").Value = 11.21653887140452
$ws.Range("E38").Value = 22.41808675646531
$ws.Range("C39").Value = 0.421655805130472
$ws.Range("E39").Value = -6.821105596638954
$ws.Range("E40").Value = -0.5376914776811237
$ws.Range("C41").Value = 5.042810166847067
$ws.Range("C43").Value = -1.548915741813695
$ws.Range("E43").Value = -3.955662492975198
$ws.Range("E44").Value = -4.829433539906869
$ws.Range("C45").Value = -0.5369231962162102
$ws.Range("E45").Value = 0.09950561885605502
$ws.Range("C46").Value = -0.9008525709169657
$ws.Range("E46").Value = -0.9756765446554017
$ws.Range("E47").Value = 0.9503229429644433
$ws.Range("C48").Value = -0.244366674180263
$ws.Range("E48").Value = -1.64927836088965
$ws.Range("E49").Value = -0.7585430378855618
$ws.Range("C50").Value = 0.2738544794132602
$ws.Range("C51").Value = 0.3683024421824888
$ws.Range("E51").Value = 0.3338002926567718
$ws.Range("E53").Value = 0.3540813801726106
